$wb = $excel.ActiveWorkbook

# --- Sheet "gof": update Deviance, AIC, BIC numeric values ---
$gof = $wb.Worksheets.Item("gof")
$gof.Range("D2").Value = 15637
$gof.Range("F2").Value = 15673
$gof.Range("G2").Value = 15775
$gof.Range("D3").Value = 15628
$gof.Range("F3").Value = 15694
$gof.Range("G3").Value = 15880
$gof.Range("D4").Value = 15574
$gof.Range("F4").Value = 15612
$gof.Range("G4").Value = 15719
$gof.Range("D5").Value = 15530
$gof.Range("F5").Value = 15628
$gof.Range("G5").Value = 15904

# --- Sheet "estimates": update DIF/estimate value cells (col B-E, rows 2-19) ---
$est = $wb.Worksheets.Item("estimates")
$est.Range("B2").Value = "-0.089 (-0.084)"
$est.Range("C2").Value = " 0.004 ( 0.004)"
$est.Range("D2").Value = " 0.005 ( 0.005)"
$est.Range("E2").Value = " 0.002 ( 0.002)"
$est.Range("B3").Value = "-0.082 (-0.077)"
$est.Range("C3").Value = " 0.158 ( 0.153)"
$est.Range("D3").Value = " 0.170 ( 0.164)"
$est.Range("E3").Value = " 0.012 ( 0.012)"
$est.Range("B4").Value = "-0.107 (-0.101)"
$est.Range("C4").Value = " 0.075 ( 0.073)"
$est.Range("D4").Value = "-0.129 (-0.125)"
$est.Range("E4").Value = "-0.205 (-0.198)"
$est.Range("B5").Value = " 0.050 ( 0.047)"
$est.Range("C5").Value = " 0.169 ( 0.163)"
$est.Range("D5").Value = "-0.290 (-0.281)"
$est.Range("E5").Value = "-0.459 (-0.444)"
$est.Range("B6").Value = "-0.022 (-0.021)"
$est.Range("C6").Value = "-0.207 (-0.200)"
$est.Range("D6").Value = "-0.351 (-0.340)"
$est.Range("E6").Value = "-0.144 (-0.139)"
$est.Range("B7").Value = "-0.042 (-0.040)"
$est.Range("C7").Value = " 0.002 ( 0.002)"
$est.Range("D7").Value = " 0.426 ( 0.412)"
$est.Range("E7").Value = " 0.425 ( 0.411)"
$est.Range("B8").Value = " 0.123 ( 0.116)"
$est.Range("C8").Value = " 0.020 ( 0.019)"
$est.Range("D8").Value = "-0.101 (-0.098)"
$est.Range("E8").Value = "-0.121 (-0.117)"
$est.Range("B9").Value = "-0.168 (-0.159)"
$est.Range("C9").Value = " 0.066 ( 0.064)"
$est.Range("D9").Value = " 0.658 ( 0.637)"
$est.Range("E9").Value = " 0.591 ( 0.572)"
$est.Range("B10").Value = "-0.044 (-0.042)"
$est.Range("C10").Value = " 0.022 ( 0.021)"
$est.Range("D10").Value = " 0.414 ( 0.401)"
$est.Range("E10").Value = " 0.392 ( 0.379)"
$est.Range("B11").Value = "-0.182 (-0.172)"
$est.Range("C11").Value = "-0.138 (-0.134)"
$est.Range("D11").Value = "-0.365 (-0.353)"
$est.Range("E11").Value = "-0.227 (-0.220)"
$est.Range("B12").Value = " 0.125 ( 0.118)"
$est.Range("C12").Value = " 0.305 ( 0.295)"
$est.Range("D12").Value = " 0.669 ( 0.647)"
$est.Range("E12").Value = " 0.364 ( 0.352)"
$est.Range("B13").Value = "-0.173 (-0.163)"
$est.Range("C13").Value = "-0.083 (-0.080)"
$est.Range("D13").Value = " 0.244 ( 0.236)"
$est.Range("E13").Value = " 0.327 ( 0.316)"
$est.Range("B14").Value = " 0.181 ( 0.171)"
$est.Range("C14").Value = "-0.266 (-0.257)"
$est.Range("D14").Value = "-0.680 (-0.658)"
$est.Range("E14").Value = "-0.413 (-0.400)"
$est.Range("B15").Value = " 0.201 ( 0.190)"
$est.Range("C15").Value = " 0.227 ( 0.220)"
$est.Range("D15").Value = " 0.186 ( 0.180)"
$est.Range("E15").Value = "-0.042 (-0.041)"
$est.Range("B16").Value = " 0.323 ( 0.305)"
$est.Range("C16").Value = " 0.055 ( 0.053)"
$est.Range("D16").Value = "-0.159 (-0.154)"
$est.Range("E16").Value = "-0.213 (-0.206)"
$est.Range("B17").Value = " 0.094 ( 0.089)"
$est.Range("C17").Value = " 0.409 ( 0.396)"
$est.Range("D17").Value = " 0.698 ( 0.675)"
$est.Range("E17").Value = " 0.289 ( 0.280)"
$est.Range("B18").Value = "-0.413 (-0.390)"
$est.Range("C18").Value = "0.727 (0.704)"
$est.Range("D18").Value = "0.305 (0.295)"
$est.Range("E18").Value = "-0.422 (-0.408)"
$est.Range("B19").Value = "-0.400 (-0.377)"
$est.Range("C19").Value = "0.707 (0.686)"
$est.Range("D19").Value = "0.298 (0.289)"
$est.Range("E19").Value = "-0.409 (-0.397)"
